$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4788.6
$ws.Range("J62").Value = 5697.6665
$ws.Range("L62").Value = 5697.6665
$ws.Range("N62").Value = -6945.6665
$ws.Range("H65").Value = 4788.6
$ws.Range("J65").Value = 5697.6665
$ws.Range("L65").Value = 28488.3325
$ws.Range("N65").Value = -34728.3325
$ws.Range("H116").Value = 350065.6
$ws.Range("I116").Value = 834984.0600000001
$ws.Range("J116").Value = 7770.1763
$ws.Range("K116").Value = 834984.0600000001
$ws.Range("L116").Value = 7770.1763
$ws.Range("M116").Value = -831542.0600000001
$ws.Range("N116").Value = -14654.1763
$ws.Range("H141").Value = 668267.3
$ws.Range("I141").Value = 1000001
$ws.Range("J141").Value = 4800
$ws.Range("K141").Value = 3000003
$ws.Range("L141").Value = 14400
$ws.Range("M141").Value = -2994823
$ws.Range("N141").Value = -24760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6393.5776
$ws.Range("I32").Value = 5484.1797
$ws.Range("K32").Value = 5484.1797
$ws.Range("M32").Value = -5197.1797
$ws.Range("H74").Value = 8774.076999999999
$ws.Range("I74").Value = 11060.375
$ws.Range("J74").Value = 5116
$ws.Range("K74").Value = 11060.375
$ws.Range("L74").Value = 5116
$ws.Range("M74").Value = -10186.375
$ws.Range("N74").Value = -6864
$ws.Range("H77").Value = 8774.076999999999
$ws.Range("I77").Value = 11060.375
$ws.Range("J77").Value = 5116
$ws.Range("K77").Value = 55301.875
$ws.Range("L77").Value = 25580
$ws.Range("M77").Value = -50933.875
$ws.Range("N77").Value = -34316

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 194.87878
$ws.Range("J80").Value = 207.52174
$ws.Range("L80").Value = 207.52174
$ws.Range("N80").Value = -2203.52174
$ws.Range("H83").Value = 194.87878
$ws.Range("J83").Value = 207.52174
$ws.Range("L83").Value = 1037.6087
$ws.Range("N83").Value = -11021.6087
$ws.Range("H134").Value = 3396.2
$ws.Range("I134").Value = 2637.4285
$ws.Range("J134").Value = 5166.6665
$ws.Range("K134").Value = 7912.2855
$ws.Range("L134").Value = 15499.9995
$ws.Range("M134").Value = -5377.2855
$ws.Range("N134").Value = -20569.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 52662.375
$ws.Range("J68").Value = 52662.375
$ws.Range("L68").Value = 52662.375
$ws.Range("N68").Value = -54160.375
$ws.Range("H71").Value = 52662.375
$ws.Range("J71").Value = 52662.375
$ws.Range("L71").Value = 157987.125
$ws.Range("N71").Value = -165475.125
$ws.Range("H99").Value = 11768976
$ws.Range("I99").Value = 16668716
$ws.Range("J99").Value = 9600
$ws.Range("K99").Value = 16668716
$ws.Range("L99").Value = 9600
$ws.Range("M99").Value = -16667218
$ws.Range("N99").Value = -12596
$ws.Range("H126").Value = 11768976
$ws.Range("I126").Value = 16668716
$ws.Range("J126").Value = 9600
$ws.Range("K126").Value = 50006148
$ws.Range("L126").Value = 28800
$ws.Range("M126").Value = -50003678
$ws.Range("N126").Value = -33740
$ws.Range("H132").Value = 4348.4546
$ws.Range("I132").Value = 4429.636
$ws.Range("J132").Value = 4267.273
$ws.Range("K132").Value = 13288.908
$ws.Range("L132").Value = 12801.819
$ws.Range("M132").Value = -10758.908
$ws.Range("N132").Value = -17861.819
$ws.Range("H134").Value = 2695.2307
$ws.Range("I134").Value = 1337.5555
$ws.Range("J134").Value = 5750
$ws.Range("K134").Value = 4012.6665
$ws.Range("L134").Value = 17250
$ws.Range("M134").Value = -1477.6665
$ws.Range("N134").Value = -22320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 100472.1
$ws.Range("I107").Value = 430
$ws.Range("J107").Value = 200514.2
$ws.Range("K107").Value = 1290
$ws.Range("L107").Value = 601542.6000000001
$ws.Range("M107").Value = 630
$ws.Range("N107").Value = -605382.6000000001
$ws.Range("H132").Value = 2364.5881

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 34970
$ws.Range("J35").Value = 34970
$ws.Range("L35").Value = 34970
$ws.Range("N35").Value = -35566
$ws.Range("H134").Value = 34901.535
$ws.Range("J134").Value = 34901.535
$ws.Range("L134").Value = 104704.605
$ws.Range("N134").Value = -109774.605

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 4673000.5
$ws.Range("J11").Value = 4673000.5
$ws.Range("L11").Value = 4673000.5
$ws.Range("N11").Value = -4673280.5
$ws.Range("H22").Value = 68859
$ws.Range("I22").Value = 168213.5
$ws.Range("J22").Value = 2622.6667
$ws.Range("K22").Value = 168213.5
$ws.Range("L22").Value = 2622.6667
$ws.Range("M22").Value = -167918.5
$ws.Range("N22").Value = -3212.6667
$ws.Range("H27").Value = 68859
$ws.Range("I27").Value = 168213.5
$ws.Range("J27").Value = 2622.6667
$ws.Range("K27").Value = 168213.5
$ws.Range("L27").Value = 2622.6667
$ws.Range("M27").Value = -168106.5
$ws.Range("N27").Value = -2836.6667
$ws.Range("H61").Value = 2070.6316
$ws.Range("I61").Value = 2228.9092
$ws.Range("J61").Value = 1853
$ws.Range("K61").Value = 2228.9092
$ws.Range("L61").Value = 1853
$ws.Range("M61").Value = -2026.9092
$ws.Range("N61").Value = -2257
$ws.Range("H68").Value = 649.03
$ws.Range("I68").Value = 649.03
$ws.Range("K68").Value = 649.03
$ws.Range("M68").Value = 99.97000000000003
$ws.Range("H71").Value = 649.03
$ws.Range("I71").Value = 649.03
$ws.Range("K71").Value = 3245.15
$ws.Range("M71").Value = 498.8500000000004
$ws.Range("H113").Value = 2070.6316
$ws.Range("I113").Value = 2228.9092
$ws.Range("J113").Value = 1853
$ws.Range("K113").Value = 2228.9092
$ws.Range("L113").Value = 1853
$ws.Range("M113").Value = -58.90920000000006
$ws.Range("N113").Value = -6193

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 27900
$ws.Range("J64").Value = 27900
$ws.Range("L64").Value = 27900
$ws.Range("N64").Value = -28396
$ws.Range("H67").Value = 27900
$ws.Range("J67").Value = 27900
$ws.Range("L67").Value = 27900
$ws.Range("N67").Value = -29616
$ws.Range("H136").Value = 7155
$ws.Range("I136").Value = 3802.5
$ws.Range("J136").Value = 10028.571
$ws.Range("K136").Value = 11407.5
$ws.Range("L136").Value = 30085.713
$ws.Range("M136").Value = -8857.5
$ws.Range("N136").Value = -35185.713
